# "commit view for account. only index and list"
# The sheet is trimmed down to just the index/list columns:
#   - the id column header becomes "stt"
#   - the trailing "Description" / "created by" / "created on" columns
#     (K:M) are removed entirely, along with their data
#   - the now-empty trailing row (row 9) is removed
#   - selection is left on the newly-emptied column K, matching the
#     state Excel leaves behind right after deleting K:M

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last three columns (Description, created by, created on).
$ws.Range("K1:M1").EntireColumn.Delete() | Out-Null

# Drop the now-trailing blank row.
$ws.Rows("9:9").Delete() | Out-Null

# Rename the id header to "stt".
$ws.Range("A1").Value2 = "stt"

# Leave the selection on column K (now blank), as Excel does after
# deleting the K:M columns.
$ws.Columns("K:K").Select() | Out-Null
